# Updated cryptos list on Sat Oct 19 22:39:32 UTC 2024 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for the crypto table, and
# swaps the Polkadot/Uniswap rows (20/21) to reflect their new rank order.
#
# Price cells are written as literal text in the source data (e.g.
# "68.330.22", "597.79"). Plain "$range.Value = '597.79'" would let Excel's
# type inference turn that into the number 597.79 (dropping the trailing
# zero / exact text). Forcing NumberFormat to Text ("@") first for the
# cells whose new value is unambiguously numeric keeps them stored as text,
# matching the original inline-string content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.330.22'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '2.649.32'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.79'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.05'
$ws.Range('E6').Value = '  +2.53%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('E9').Value = '  +3.02%  '
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.352'
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.12'
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').Value = '3.127.60'
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('D16').Value = '68.208.34'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '2.648.17'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.41'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '364.32'
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.33'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.43'
$ws.Range('E21').Value = '  +3.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.81'
$ws.Range('E22').Value = '  -0.55%  '
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.30'
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.74'
$ws.Range('E26').Value = '  -2.69%  '
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '558.41'
$ws.Range('E30').Value = '  -2.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.02'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +1.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.85'
$ws.Range('E37').Value = '  +3.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '159.91'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.372'
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.88'
$ws.Range('E40').Value = '  -2.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.34'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('E42').Value = '  +3.60%  '
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '158.62'
$ws.Range('E45').Value = '  +1.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.76'
$ws.Range('E46').Value = '  +0.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.14'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.567'
$ws.Range('E51').Value = '  +1.13%  '
